$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing used range (rows 1-7, cols A-C) before writing new data
$ws.Range("A1:C7").Clear()

# Force text format so values are stored as strings (inlineStr), not numbers
$ws.Range("A1:D4").NumberFormat = "@"

# New header row
$ws.Range("A1").Value = "NrComClient"
$ws.Range("B1").Value = "PretPropusClient"
$ws.Range("C1").Value = "CodInitialComClient"
$ws.Range("D1").Value = "CantInitialaComClient"

# Row 2
$ws.Range("A2").Value = "20250628104637"
$ws.Range("B2").Value = "17.22"
$ws.Range("C2").Value = "LR067042"
$ws.Range("D2").Value = "1"

# Row 3
$ws.Range("A3").Value = "20250628124037"
$ws.Range("B3").Value = "34"
$ws.Range("C3").Value = "30735186"
$ws.Range("D3").Value = "1"

# Row 4
$ws.Range("A4").Value = "20250628123337"
$ws.Range("B4").Value = "59.29"
$ws.Range("C4").Value = "281132S0"
$ws.Range("D4").Value = "1"
